# feat: add 2022-Q4 data
#
# The workbook currently has two sheets:
#   总计      - summary sheet with one row of totals per quarter
#   2022-Q3   - fund-holding detail sheet for 2022-Q3
#
# We are adding a new quarter (2022-Q4). The new quarter's data takes over the
# existing "2022-Q3" detail sheet's slot (renamed to "2022-Q4" with updated
# figures), while a fresh copy of the original "2022-Q3" sheet (with its
# original data untouched) is inserted right after it, so the historical
# 2022-Q3 snapshot is preserved as its own tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q3" detail sheet so its current data is
#    preserved on its own tab, placed immediately after the source sheet.
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($null, $q3Sheet)

# The duplicate is placed right after $q3Sheet, i.e. it is now the last sheet.
$q3Copy = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------
# 2) Turn the original "2022-Q3" sheet into the new "2022-Q4" sheet by
#    renaming it and overwriting its fund-holding figures with the new
#    quarter's values (fund code/name stay the same). Rename it before the
#    duplicate so the two sheets never momentarily share a name.
# ---------------------------------------------------------------------
$q3Sheet.Name = "2022-Q4"
$q3Copy.Name = "2022-Q3"

$q3Sheet.Range("D2").Value = 9.01
$q3Sheet.Range("E2").Value = 94.55
$q3Sheet.Range("F2").Value = 4.27
$q3Sheet.Range("G2").Value = 0.3847

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: the row that used to describe
#    2022-Q3 now describes 2022-Q4 with its new totals, and a new row is
#    appended below it preserving the old 2022-Q3 totals.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q4"
$total.Range("D2").Value = 0.38

# Copy A2's formatting (style s="2") down into A3, then set the new row's
# values (historical 2022-Q3 totals).
$total.Range("A2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.34
